# Historique.xlsx - add the "Chargement des niveaux depuis un fichier" task
# to the task table on the first sheet (row 16), with its time estimate,
# the TP it was realised in, and its progress status.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("A16").Value = "Chargement des niveaux depuis un fichier"
$ws.Range("B16").Value = "1h30"
$ws.Range("C16").Value = "TP3"
$ws.Range("F16").Value = "OK"

# Scroll the view down a bit and leave the newly-edited cell selected,
# mirroring what the author would see after typing the row in.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F16").Select()
